$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("113_2")

# Update the End Date row (row 5) values from 2/2/2015 (42037) to 1/2/2015 (42006)
$ws.Range("B5:D5").Value = "1/2/2015"

$wb.Save()
